$wb = $excel.ActiveWorkbook

# The localization run for this file moved on from handoff into active
# translation, so the cached "Status" (and per-locale rollup on the
# Overview sheet) needs to read "In Translation" instead of
# "Ready for handoff".
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# The status text is shorter now, so the status/locale columns get
# narrowed to fit the new value.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
